$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "fgfd"
$ws.Range("K4").Value = "fgdsgs"
$ws.Range("I8").Value = "sgfs"

$ws.Range("I8").Select()
